# Scheduled market-data refresh: update Leve profit calculations (currentAveragePrice*,
# LevePrice*, LeveProfit* columns) per-job across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 74
$ws.Range("H74").Value = 3500
$ws.Range("J74").Value = 3500
$ws.Range("L74").Value = 3500
$ws.Range("N74").Value = -5372
# Row 77
$ws.Range("H77").Value = 3500
$ws.Range("J77").Value = 3500
$ws.Range("L77").Value = 17500
$ws.Range("N77").Value = -26860
# Row 130
$ws.Range("H130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("L130").Value = 0
$ws.Range("N130").ClearContents()
$ws = $wb.Worksheets.Item("ARM")
# Row 45
$ws.Range("H45").Value = 1798.6666
$ws.Range("J45").Value = 2200
$ws.Range("L45").Value = 2200
$ws.Range("N45").Value = -2954
# Row 61
$ws.Range("H61").Value = 1878.1
$ws.Range("I61").Value = 1878.1
$ws.Range("K61").Value = 1878.1
$ws.Range("M61").Value = -1666.1
# Row 102
$ws.Range("H102").Value = 1101.5333
$ws.Range("I102").Value = 937.3570999999999
$ws.Range("K102").Value = 937.3570999999999
$ws.Range("M102").Value = 684.6429000000001
# Row 110
$ws.Range("H110").Value = 949.5
$ws.Range("I110").Value = 943.8889
$ws.Range("K110").Value = 943.8889
$ws.Range("M110").Value = 1101.1111
# Row 122
$ws.Range("H122").Value = 4499.5
$ws.Range("I122").Value = 4499.5
$ws.Range("K122").Value = 13498.5
$ws.Range("M122").Value = -11048.5
# Row 136
$ws.Range("H136").Value = 1878.1
$ws.Range("I136").Value = 1878.1
$ws.Range("K136").Value = 5634.299999999999
$ws.Range("M136").Value = -3084.299999999999
$ws = $wb.Worksheets.Item("BSM")
# Row 86
$ws.Range("H86").Value = 7221.4443
$ws.Range("I86").Value = 7333
$ws.Range("J86").Value = 7165.6665
$ws.Range("K86").Value = 7333
$ws.Range("L86").Value = 7165.6665
$ws.Range("M86").Value = -6210
$ws.Range("N86").Value = -9411.666499999999
# Row 89
$ws.Range("H89").Value = 7221.4443
$ws.Range("I89").Value = 7333
$ws.Range("J89").Value = 7165.6665
$ws.Range("K89").Value = 36665
$ws.Range("L89").Value = 35828.3325
$ws.Range("M89").Value = -31049
$ws.Range("N89").Value = -47060.3325
# Row 105
$ws.Range("H105").Value = 3172.1333
$ws.Range("I105").Value = 3028.5
$ws.Range("K105").Value = 3028.5
$ws.Range("M105").Value = -1281.5
# Row 107
$ws.Range("H107").Value = 1073.0625
$ws.Range("I107").Value = 859.53845
$ws.Range("J107").Value = 1998.3334
$ws.Range("K107").Value = 859.53845
$ws.Range("L107").Value = 1998.3334
$ws.Range("M107").Value = 1060.46155
$ws.Range("N107").Value = -5838.3334
$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 977.4545000000001
$ws.Range("J16").Value = 442.5
$ws.Range("L16").Value = 442.5
$ws.Range("N16").Value = -1016.5
# Row 58
$ws.Range("H58").Value = 4741.9
$ws.Range("I58").Value = 3828.0625
$ws.Range("K58").Value = 3828.0625
$ws.Range("M58").Value = -3625.0625
# Row 86
$ws.Range("H86").Value = 9490.5
$ws.Range("I86").Value = 8722
$ws.Range("K86").Value = 8722
$ws.Range("M86").Value = -7599
# Row 89
$ws.Range("H89").Value = 9490.5
$ws.Range("I89").Value = 8722
$ws.Range("K89").Value = 43610
$ws.Range("M89").Value = -37994
# Row 94
$ws.Range("H94").Value = 1559.8572
$ws.Range("I94").Value = 1536.6666
$ws.Range("K94").Value = 1536.6666
$ws.Range("M94").Value = -1085.6666
# Row 113
$ws.Range("H113").Value = 977.4545000000001
$ws.Range("J113").Value = 442.5
$ws.Range("L113").Value = 442.5
$ws.Range("N113").Value = -4782.5
# Row 136
$ws.Range("H136").Value = 4741.9
$ws.Range("I136").Value = 3828.0625
$ws.Range("K136").Value = 11484.1875
$ws.Range("M136").Value = -8934.1875
$ws = $wb.Worksheets.Item("CUL")
# Row 51
$ws.Range("H51").Value = 898.6667
$ws.Range("I51").Value = 999
$ws.Range("J51").Value = 848.5
$ws.Range("K51").Value = 2997
$ws.Range("L51").Value = 2545.5
$ws.Range("M51").Value = -2537
$ws.Range("N51").Value = -3465.5
# Row 60
$ws.Range("H60").Value = 238.4
$ws.Range("I60").Value = 273
$ws.Range("J60").Value = 100
$ws.Range("K60").Value = 819
$ws.Range("L60").Value = 300
$ws.Range("M60").Value = -568
$ws.Range("N60").Value = -802
$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Range("H70").Value = 111116110
$ws.Range("I70").Value = 166669170
$ws.Range("K70").Value = 166669170
$ws.Range("M70").Value = -166668900
# Row 73
$ws.Range("H73").Value = 111116110
$ws.Range("I73").Value = 166669170
$ws.Range("K73").Value = 166669170
$ws.Range("M73").Value = -166668234
# Row 113
$ws.Range("H113").Value = 2473.6667
$ws.Range("J113").Value = 2210
$ws.Range("L113").Value = 2210
$ws.Range("N113").Value = -6550
$ws = $wb.Worksheets.Item("LTW")
# Row 61
$ws.Range("H61").Value = 3099.8
$ws.Range("I61").Value = 2499.75
$ws.Range("J61").Value = 5500
$ws.Range("K61").Value = 2499.75
$ws.Range("L61").Value = 5500
$ws.Range("M61").Value = -2297.75
$ws.Range("N61").Value = -5904
# Row 68
$ws.Range("H68").Value = 26774.75
$ws.Range("I68").Value = 1739.8
$ws.Range("K68").Value = 1739.8
$ws.Range("M68").Value = -990.8
# Row 71
$ws.Range("H71").Value = 26774.75
$ws.Range("I71").Value = 1739.8
$ws.Range("K71").Value = 8699
$ws.Range("M71").Value = -4955
# Row 82
$ws.Range("H82").Value = 1161
$ws.Range("I82").Value = 1395.6666
$ws.Range("J82").Value = 926.3333
$ws.Range("K82").Value = 1395.6666
$ws.Range("L82").Value = 926.3333
$ws.Range("M82").Value = -1034.6666
$ws.Range("N82").Value = -1648.3333
# Row 85
$ws.Range("H85").Value = 1161
$ws.Range("I85").Value = 1395.6666
$ws.Range("J85").Value = 926.3333
$ws.Range("K85").Value = 1395.6666
$ws.Range("L85").Value = 926.3333
$ws.Range("M85").Value = -147.6666
$ws.Range("N85").Value = -3422.3333
# Row 113
$ws.Range("H113").Value = 3099.8
$ws.Range("I113").Value = 2499.75
$ws.Range("J113").Value = 5500
$ws.Range("K113").Value = 2499.75
$ws.Range("L113").Value = 5500
$ws.Range("M113").Value = -329.75
$ws.Range("N113").Value = -9840
# Row 136
$ws.Range("H136").Value = 2304.3333
$ws.Range("I136").Value = 1886.8
$ws.Range("K136").Value = 5660.4
$ws.Range("M136").Value = -3110.4
$ws = $wb.Worksheets.Item("WVR")
# Row 55
$ws.Range("H55").Value = 25250
$ws.Range("J55").Value = 25250
$ws.Range("L55").Value = 25250
$ws.Range("N55").Value = -25804
# Row 81
$ws.Range("H81").Value = 2761.25
$ws.Range("I81").Value = 2761.25
$ws.Range("K81").Value = 5522.5
$ws.Range("M81").Value = -4461.5
# Row 84
$ws.Range("H84").Value = 2761.25
$ws.Range("I84").Value = 2761.25
$ws.Range("K84").Value = 27612.5
$ws.Range("M84").Value = -22308.5
# Row 132
$ws.Range("H132").Value = 1277.1666
$ws.Range("I132").Value = 1262.125
$ws.Range("K132").Value = 3786.375
$ws.Range("M132").Value = -1256.375
